$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp string (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 8 de Agosto de 2020 a las 07:57"

# --- Swap the country labels for Timor Oriental / Santa Lucia rows ---
# (their numeric stats are identical, only the shared-string order/labels swap)
$ws.Range("A202").Value = "Santa Lucia"
$ws.Range("A203").Value = "Timor Oriental"

# --- Pakistan (row 17) ---
$ws.Range("B17").Value = 283487
$ws.Range("C17").Value = 842
$ws.Range("D17").Value = 259604
$ws.Range("E17").Value = 17815
$ws.Range("G17").Value = 16
$ws.Range("H17").Value = 6068

# --- Kirguistan (row 56) ---
$ws.Range("B56").Value = 39571
$ws.Range("C56").Value = 409
$ws.Range("D56").Value = 31062
$ws.Range("E56").Value = 7050
$ws.Range("G56").Value = 8
$ws.Range("H56").Value = 1459

# --- Uzbekistan (row 62) ---
$ws.Range("B62").Value = 29225
$ws.Range("C62").Value = 168
$ws.Range("D62").Value = 20380
$ws.Range("E62").Value = 8662
$ws.Range("G62").Value = 2
$ws.Range("H62").Value = 183

# --- El Salvador (row 73) ---
$ws.Range("D73").Value = 9392
$ws.Range("E73").Value = 9616
$ws.Range("G73").Value = 16
$ws.Range("H73").Value = 536

# --- Tailandia (row 115) ---
$ws.Range("B115").Value = 3348
$ws.Range("C115").Value = 3
$ws.Range("D115").Value = 3150
$ws.Range("E115").Value = 140
